# Swap the contents of column E ("codeforiati:category-name") and
# column F ("codeforiati:group-code"), header row included, on the
# active worksheet.
#
# A plain Value/Value2 round-trip would make Excel re-interpret the
# numeric-looking "group-code" text (e.g. "110") as a real number once it
# lands in column E, which would change the cell's stored type. Routing
# the swap through Copy + PasteSpecial(xlPasteValues) instead copies the
# cell's existing value verbatim (keeping text as text) and does not
# touch the cell's style, so no new number-format style gets introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

$xlPasteValues = -4163

$colE = $ws.Range("E1:E$lastRow")
$colF = $ws.Range("F1:F$lastRow")
$scratch = $ws.Range("I1:I$lastRow")

$colE.Copy()
$scratch.PasteSpecial($xlPasteValues)

$colF.Copy()
$colE.PasteSpecial($xlPasteValues)

$scratch.Copy()
$colF.PasteSpecial($xlPasteValues)

$scratch.ClearContents()
$excel.CutCopyMode = $false
